# Actualización automática del mapa (2025-08-30 08:22:33)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update observation text on row 41 (column H)
$ws.Range("H41").Value = "Pendiente de traspaso nodo entro tambien como 7100"

# 2) Append two new incident rows (83 and 84) carrying the same column
#    layout as the rest of the sheet: A Caso, B F. De Reclamo, C Direccion,
#    D Comuna, E OT, F Proveedor Asignado, G Estado, H Observaciones,
#    I Attachments, J Tipo de tarea, K Equipo, L Tipo de Elemento,
#    M Coordenada_X, N Coordenada_Y, O Operacion, P Zona.

function Set-TextCell($sheet, $addr, $text) {
    # Force text storage (not auto-converted to a number/date by Excel),
    # matching how every other row on this sheet stores these columns.
    # The leading apostrophe forces Excel to treat the entry as text even
    # when it looks numeric/date-like; resetting the style back to Normal
    # afterwards keeps the cell's formatting identical to its neighbours
    # (no stray quote-prefix / number-format style gets left behind).
    $cell = $sheet.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 83
Set-TextCell $ws "A83" "4698"
Set-TextCell $ws "B83" "8/29/2025"
Set-TextCell $ws "C83" "REPETTO, NICOLAS, DR. 93"
Set-TextCell $ws "D83" "6"
Set-TextCell $ws "E83" "Pendiente ADM"
Set-TextCell $ws "F83" "PEBCOM"
Set-TextCell $ws "G83" "Pendiente de Traspaso PROPIO"
Set-TextCell $ws "H83" "traspasar nodo a columna nueva y pasar a retirar entro directamente con la nueva al lado"
$ws.Range("I83").Value = 1
Set-TextCell $ws "J83" "Cambio"
Set-TextCell $ws "K83" "Nodo Teco"
Set-TextCell $ws "L83" "Terminal"
$ws.Range("M83").Value = -58.443232
$ws.Range("N83").Value = -34.620007
Set-TextCell $ws "O83" "Boedo"
Set-TextCell $ws "P83" "Capital Sur"

# Row 84
Set-TextCell $ws "A84" "7102"
Set-TextCell $ws "B84" "8/30/2025"
Set-TextCell $ws "C84" "AMBERES 995"
Set-TextCell $ws "D84" "6"
Set-TextCell $ws "E84" "809309598"
Set-TextCell $ws "F84" "PEBCOM"
Set-TextCell $ws "G84" "Pendiente"
Set-TextCell $ws "H84" "Picada"
$ws.Range("I84").Value = 1
Set-TextCell $ws "J84" "Cambio"
Set-TextCell $ws "K84" "Sin equipos"
Set-TextCell $ws "L84" "Pasante"
$ws.Range("M84").Value = -58.453382
$ws.Range("N84").Value = -34.612707
Set-TextCell $ws "O84" "Boedo"
Set-TextCell $ws "P84" "Capital Sur"

Write-Output "Done"
